{"js": "// The document originally contains a single paragraph built from several\n// runs (plus proofErr markers) that together read:\n//   \"Hey..I work with WagMob and I love working with Git. It's a simple,\n//    fast, and superb version control system.\"\n//\n// The edit:\n//   1. Inserts a brand-new paragraph BEFORE that paragraph containing the\n//      exact same sentence (now as one clean run).\n//   2. Replaces the content of the original (now second) paragraph with a\n//      new sentence: \"Git is awesome ,it improves my productivity.\"\n\nconst body = context.document.body;\n\n// Load the existing paragraphs so we can anchor the insertion on the\n// first (and, at this point, only) paragraph in the body.\nlet paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\n\n// Insert the new paragraph with the original greeting text right before\n// the existing paragraph.\nfirstParagraph.insertParagraph(\n  \"Hey..I work with WagMob and I love working with Git. It's a simple, fast, and superb version control system.\",\n  Word.InsertLocation.before\n);\nawait context.sync();\n\n// Re-query paragraphs fresh from the body: index 0 is the newly inserted\n// paragraph, index 1 is the original paragraph whose text we still need\n// to replace.\nparagraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst secondParagraph = paragraphs.items[1];\nsecondParagraph.insertText(\n  \"Git is awesome ,it improves my productivity.\",\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n", "ps1": "# The document originally contains a single paragraph built from several\n# runs (plus proofErr markers) that together read:\n#   \"Hey..I work with WagMob and I love working with Git. It's a simple,\n#    fast, and superb version control system.\"\n#\n# The edit:\n#   1. Inserts a brand-new paragraph BEFORE that paragraph containing the\n#      exact same sentence (now as one clean run).\n#   2. Replaces the content of the original (now second) paragraph with a\n#      new sentence: \"Git is awesome ,it improves my productivity.\"\n\n$d = $word.ActiveDocument\n\n# Anchor on the (currently only) first paragraph and insert a brand-new\n# paragraph immediately before it.\n$firstParagraph = $d.Paragraphs(1)\n$firstParagraph.Range.InsertParagraphBefore()\n\n# The freshly-inserted paragraph is now paragraph 1; fill it with the\n# original greeting sentence.\n$d.Paragraphs(1).Range.Text = \"Hey..I work with WagMob and I love working with Git. It's a simple, fast, and superb version control system.\"\n\n# Paragraph 2 is the original paragraph (still holding the old runs /\n# proofErr markers). Replace its text - excluding the trailing paragraph\n# mark - with the new sentence.\n$secondRange = $d.Paragraphs(2).Range\n$secondRange.End = $secondRange.End - 1\n$secondRange.Text = \"Git is awesome ,it improves my productivity.\"\n"}
